# PV-261: Correct PV-Test-03 so it has a valid sheet name.
#
# The first worksheet (rId1 / sheetId 1) was incorrectly named
# "PV-Test-01" - rename it to "PV-Test-03" so the file-read logic (which
# picks the currently active sheet) finds a worksheet whose name matches
# the file name. Also make this corrected sheet the active/selected tab
# instead of the "Dummy" sheet.

$wb = $excel.ActiveWorkbook

# The mis-named worksheet is the first sheet in the workbook.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "PV-Test-03"

# Make this sheet the active one (moves tabSelected from "Dummy" to this
# sheet and clears the activeTab override at the workbook level).
$ws.Activate()
